$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Cells that change type/style (text<->number): fix formatting via PasteSpecial(xlPasteFormats) first ---
$ws.Range("N14").Copy()
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("F15").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("N14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("N14").Copy()
$ws.Range("L28").PasteSpecial(-4122)
$ws.Range("N14").Copy()
$ws.Range("L29").PasteSpecial(-4122)

# C26 and F30 go from NUMBER to TEXT ("0"): copy both format AND value (as text) from
# D14, which already holds the literal text "0" with the desired (General/style-14) format.
$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("F30").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Set values (numbers and the two text cells) ---
$ws.Range("L14").Value = -100
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = 0
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 166.666666666667
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 161.538461538462
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = 116.666666666667
$ws.Range("L16").Value = 136.363636363636
$ws.Range("M16").Value = -18.75
$ws.Range("N16").Value = -70.114942528735
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -45.454545454545
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -28.571428571428
$ws.Range("I17").Value = 21
$ws.Range("J17").Value = 27
$ws.Range("K17").Value = -22.222222222222
$ws.Range("L17").Value = 16.666666666666
$ws.Range("M17").Value = 10.526315789473
$ws.Range("N17").Value = -50
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -13.636363636363
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 21
$ws.Range("K18").Value = -23.809523809523
$ws.Range("L18").Value = -5.882352941176
$ws.Range("M18").Value = -20
$ws.Range("N18").Value = -79.746835443038
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = -68
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = -46.376811594202
$ws.Range("I19").Value = 27
$ws.Range("J19").Value = 56
$ws.Range("K19").Value = -51.785714285714
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 80
$ws.Range("N19").Value = -10
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 63.636363636363
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = 20
$ws.Range("L20").Value = 20
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -81.538461538461
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -36
$ws.Range("F21").Value = 135
$ws.Range("G21").Value = 154
$ws.Range("H21").Value = -12.337662337662
$ws.Range("I21").Value = 104
$ws.Range("J21").Value = 130
$ws.Range("K21").Value = -20
$ws.Range("L21").Value = 20.930232558139
$ws.Range("M21").Value = 5.050505050505
$ws.Range("N21").Value = -66.451612903225
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = 300
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = -20
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -5.882352941176
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = -9.210526315789
$ws.Range("I24").Value = 43
$ws.Range("J24").Value = 55
$ws.Range("K24").Value = -21.818181818181
$ws.Range("L24").Value = -18.867924528301
$ws.Range("M24").Value = 16.216216216216
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 55.555555555555
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 37.142857142857
$ws.Range("I25").Value = 30
$ws.Range("J25").Value = 28
$ws.Range("K25").Value = 7.142857142857
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -11.764705882352
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = -25
$ws.Range("L26").Value = 0
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 100
$ws.Range("L28").Value = -100
$ws.Range("L29").Value = -100
